$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: merge the runs around "proportion of wealth in stock i"
# so the proofErr spell-check wrapper around the lone "i" is gone and
# everything becomes a single run (text content is unchanged).
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    " proportion of wealth in stock i" + [char]0x201D + ". I" + [char]0x2019 + "m assuming",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " proportion of wealth in stock i" + [char]0x201D + ". I" + [char]0x2019 + "m assuming",
    2) | Out-Null

# -----------------------------------------------------------------
# Change 2: insert a new sub-bullet "How long is our in sample
# period" right after "How long after are we going to check our
# progress" (same ilvl = 1). The Word "_GoBack" bookmark should end
# up sitting right after the new text, exactly as it would after a
# user finishes typing it.
# -----------------------------------------------------------------
$progressPara = $d.Paragraphs(8)
$progressPara.Range.InsertParagraphAfter() | Out-Null
$newPara1 = $d.Paragraphs(9)
# Type with a trailing sentinel character first -- writing straight
# at the paragraph-end position is unreliable, so we park the
# bookmark one character to the left of the final paragraph mark and
# then trim the sentinel back off.
$newPara1.Range.Text = "How long is our in sample period#"
$r1 = $d.Paragraphs(9).Range
$bmPos1 = $r1.End - 2

$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete() | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos1, $bmPos1)) | Out-Null

$sentinelRange1 = $d.Range($bmPos1, $bmPos1 + 1)
$sentinelRange1.Delete() | Out-Null

# -----------------------------------------------------------------
# Change 3: insert a new bullet "Select a few randomly from market?"
# right after "Reinvest? " (back out to ilvl = 0).
# -----------------------------------------------------------------
$reinvestPara = $d.Paragraphs(15)
$reinvestPara.Range.InsertParagraphAfter() | Out-Null
$newPara2 = $d.Paragraphs(16)
$newPara2.Range.ListFormat.ListLevelNumber = 1
$newPara2.Range.Text = "Select a few randomly from market?"

# -----------------------------------------------------------------
# Change 4: merge the runs around "and also returns etc?" so the
# proofErr spell-check wrapper around "etc" is gone (text unchanged).
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    " and also returns etc?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " and also returns etc?",
    2) | Out-Null
